$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column contains text values that look numeric (e.g.
# "332.14"). A plain .Value assignment would let Excel auto-convert
# those into real numbers (losing formatting / exact text, e.g.
# "332.14" -> 332.1399999...). Force the destination cell to Text
# format first so the literal digits are stored verbatim, matching
# how the source data was authored. Values that already contain a
# second "." (e.g. "27.575.35") or non-numeric characters (URLs,
# names, the "%" rows) are never auto-converted, so they are left
# on the sheets default General format.

$ws.Range("D5,D6,D7,D8,D9,D10,D11,D12,D14,D15,D17,D18,D19,D20,D23,D24,D26,D27,D28,D29,D30,D31,D32,D34,D35,D36,D37,D38,D40,D41,D42,D43,D44,D45,D47,D48,D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '27.575.35'
$ws.Range("E2").Value = '  -1.38%  '

# Row 3
$ws.Range("D3").Value = '1.846.28'
$ws.Range("E3").Value = '  -2.27%  '

# Row 4
$ws.Range("E4").Value = '  -1.30%  '

# Row 5
$ws.Range("D5").Value = '332.14'
$ws.Range("E5").Value = '  -1.20%  '

# Row 6
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  -1.24%  '

# Row 7
$ws.Range("D7").Value = '0.4617'
$ws.Range("E7").Value = '  -2.22%  '

# Row 8
$ws.Range("D8").Value = '0.3847'
$ws.Range("E8").Value = '  -2.50%  '

# Row 9
$ws.Range("D9").Value = '45.79'
$ws.Range("E9").Value = '  -2.37%  '

# Row 10
$ws.Range("D10").Value = '0.07907'
$ws.Range("E10").Value = '  -1.30%  '

# Row 11
$ws.Range("D11").Value = '0.9915'
$ws.Range("E11").Value = '  -2.69%  '

# Row 12
$ws.Range("D12").Value = '21.45'
$ws.Range("E12").Value = '  -1.63%  '

# Row 13
$ws.Range("D13").Value = '1.855.16'
$ws.Range("E13").Value = '  -2.45%  '

# Row 14
$ws.Range("D14").Value = '5.898'
$ws.Range("E14").Value = '  -1.66%  '

# Row 15
$ws.Range("D15").Value = '7.085'
$ws.Range("E15").Value = '  -1.43%  '

# Row 16
$ws.Range("E16").Value = '  -1.27%  '

# Row 17
$ws.Range("D17").Value = '88.56'
$ws.Range("E17").Value = '  +0.49%  '

# Row 18
$ws.Range("D18").Value = '0.06651'
$ws.Range("E18").Value = '  -2.12%  '

# Row 19
$ws.Range("D19").Value = '0.00001033'
$ws.Range("E19").Value = '  -1.80%  '

# Row 20
$ws.Range("D20").Value = '16.96'
$ws.Range("E20").Value = '  -1.18%  '

# Row 21
$ws.Range("E21").Value = '  -1.18%  '

# Row 22
$ws.Range("D22").Value = '27.583.12'
$ws.Range("E22").Value = '  -1.38%  '

# Row 23
$ws.Range("D23").Value = '5.369'
$ws.Range("E23").Value = '  -2.51%  '

# Row 24
$ws.Range("D24").Value = '10.91'
$ws.Range("E24").Value = '  -0.80%  '

# Row 25
$ws.Range("E25").Value = '  -2.37%  '

# Row 26
$ws.Range("D26").Value = '157.78'
$ws.Range("E26").Value = '  -1.16%  '

# Row 27
$ws.Range("D27").Value = '19.46'
$ws.Range("E27").Value = '  -2.85%  '

# Row 28
$ws.Range("D28").Value = '2.085'
$ws.Range("E28").Value = '  -1.18%  '

# Row 29
$ws.Range("D29").Value = '5.388'
$ws.Range("E29").Value = '  -2.19%  '

# Row 30
$ws.Range("D30").Value = '119.57'
$ws.Range("E30").Value = '  -1.70%  '

# Row 31
$ws.Range("D31").Value = '0.9713'
$ws.Range("E31").Value = '  +0.46%  '

# Row 32
$ws.Range("D32").Value = '0.09368'
$ws.Range("E32").Value = '  -2.18%  '

# Row 33
$ws.Range("E33").Value = '  -1.94%  '

# Row 34
$ws.Range("D34").Value = '5.268'
$ws.Range("E34").Value = '  -1.73%  '

# Row 35
$ws.Range("D35").Value = '1.341'
$ws.Range("E35").Value = '  -1.59%  '

# Row 36
$ws.Range("D36").Value = '0.05998'
$ws.Range("E36").Value = '  -2.02%  '

# Row 37
$ws.Range("D37").Value = '0.02216'
$ws.Range("E37").Value = '  -1.50%  '

# Row 38
$ws.Range("D38").Value = '8.281'
$ws.Range("E38").Value = '  -0.37%  '

# Row 39
$ws.Range("E39").Value = '  -2.98%  '

# Row 40
$ws.Range("D40").Value = '0.5879'
$ws.Range("E40").Value = '  -1.57%  '

# Row 41
$ws.Range("D41").Value = '0.1858'
$ws.Range("E41").Value = '  -2.78%  '

# Row 42
$ws.Range("D42").Value = '10.26'
$ws.Range("E42").Value = '  -1.36%  '

# Row 43
$ws.Range("D43").Value = '1.255'
$ws.Range("E43").Value = '  -1.63%  '

# Row 44
$ws.Range("D44").Value = '0.5571'
$ws.Range("E44").Value = '  -2.31%  '

# Row 45
$ws.Range("D45").Value = '12.07'
$ws.Range("E45").Value = '  -0.85%  '

# Row 46
$ws.Range("E46").Value = '  -2.85%  '

# Row 47
$ws.Range("D47").Value = '0.06686'
$ws.Range("E47").Value = '  -2.66%  '

# Row 48
$ws.Range("D48").Value = '110.37'
$ws.Range("E48").Value = '  -2.71%  '

# Row 49
$ws.Range("E49").Value = '  -1.86%  '

# Row 50
$ws.Range("E50").Value = '  -1.41%  '

# Row 51
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").Value = '0.00000000287'
$ws.Range("E51").Value = '  -3.27%  '
